$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 'Amity Elementary School'
$ws.Range("B5").Value = 'Amity SD 4J'
$ws.Range("C5").Value = 'Diana Sohn'
$ws.Range("D5").Value = 'Diana Sohn'
$ws.Range("E5").Value = 'diana.sohn@amity.k12.or.us'
$ws.Range("F5").Value = 'diana.sohn@amity.k12.or.us'
$ws.Range("G5").Value = '503-835-3751 ext 314'
$ws.Range("H5").Value = '503-835-3751 ext. 314'
$ws.Range("I5").Value = 4
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = 5
$ws.Range("N5").Value = 5
$ws.Range("O5").Value = 3
$ws.Range("P5").Value = 4
$ws.Range("Q5").Value = 4
$ws.Range("R5").Value = 3
$ws.Range("S5").Value = 3
$ws.Range("T5").Value = 5
$ws.Range("U5").Value = 4
$ws.Range("V5").Value = 3
$ws.Range("W5").Value = 5
$ws.Range("X5").Value = 4
$ws.Range("Y5").Value = 5
$ws.Range("AA5").Value = 'We provide a number of hands-on inquiry based activities and experiences that are integrated with our science, technology, engineering, and mathematics curricula.  For example:  building shelters, determining calories, fishing,…'
$ws.Range("AB5").Value = 'We have at least 1 high school counselor who helps with translations during ODS meetings for parents and is with the ELL student during camp to help with understanding of the content/concepts. '
$ws.Range("AC5").Value = 'Students and families are not charged money to participate. Fundraising has been done by the students to keep this program running for over 35 years.   All 5th grade children attend and have done so for over 35 years.  No one is excluded based on socio-economic status.  Many of our experiences take place in areas that our students have not had the opportunity to visit. '
$ws.Range("AD5").Value = 'Amity is a rural community, during Outdoor School, we visit a variety of terrain (forest, desert, cave, National Park, State Park, museums, nature centers) to show the differences and similarities throughout Oregon and how we depend on each other.'
$ws.Range("AE5").Value = 'x'
$ws.Range("AF5").Value = 'x'
$ws.Range("AG5").Value = 'x'
$ws.Range("AH5").Value = 'x'
$ws.Range("AI5").Value = 'x'
$ws.Range("AJ5").Value = 'x'
$ws.Range("AL5").Value = 'Enter "other" text here.'
$ws.Range("AM5").Value = 'Student journals are used when we get back to reflect upon our experiences. Journals are also used to interact with other students about what they learned. Pre- ODS we visited the Evergreen Museum and Cruickshank Farms to gain some pre knowledge about space and timber topics. Post- writing thank you letters to counselors, parents, bus drivers, and teachers. Students also wrote an ODS memory reflecting on their favorite activity from camp. Our science curriculum revolves around our Outdoor School and is integrated throughout the year.'
$ws.Range("AN5").Value = 'The chaperones are high school students grades 10-12th with good grades and approval of their teachers. High School students are responsible for a small tent group and teaching a study session six times.  They are on duty 23 hours a day for four days.  They also receive 20 hours of training prior to Outdoor School.  This past year over 40 parents/grandparents attended and helped with kitchen duties, recreation, and setup and break down of camp. One of the bus drivers is the chief of the fire department in town, who volunteers his time.  We also have at least four community members who volunteer to help with cooking.'
$ws.Range("AO5").Value = 'In science, 68% of students passed the OAKS Science test due to their experiences at ODS. Without these experiences students would not have done as well.  Students have learned how to properly write thank you notes and reflections.'
$ws.Range("AP5").Value = 'Many students volunteered hours upon hours to help prepare for ODS. It empowered students to become leaders in their community and in the classroom. Students that volunteered were able to submit hours for their Presidential Volunteer Pin. The high school counselors received over 100 hours of volunteer hours-which displays their leadership skills in running study sessions and leading their groups. '
$ws.Range("AQ5").Value = 'Because students know they have to be respectful, responsible, and safe to attend Outdoor School, classroom management and discipline problems are minimal.  Students who have behavior or missing work issues are placed on Outdoor School contracts and the majority of these students improve their behaviors.'
$ws.Range("AR5").Value = 'Students earn beads and awards in camp that transfer back to the classroom when we return.  It builds community within their class.'
$ws.Range("AS5").Value = 'Outdoor School provides students the opportunity to improve their problem solving and strategic thinking due to the open ended experiences.  For example, in survival - how to build a shelter, how to fix their leaking tent, how to take care of personal needs independently, and how to think fast when it came to a thunder and lightning storm.'
$ws.Range("AT5").Value = 'Students had opportunities to create shelters, make and present songs and skits, weave a design on a loom, build a nutritious parfait, create an obstacle course.'
$ws.Range("AU5").Value = 'Students were in a variety of groups during study sessions, tent groups, field trips away from camps, recreation, and partner work where they had to communicate with their peers in order to successfully complete tasks.'
$ws.Range("AV5").Value = 'Our students were very enthusiastic about science (space, birds, volcanic activity, animals) and social studies (High Desert Museum).  Looking through telescopes, observing animals, exploring caves resulted in excitement about learning new things.'
$ws.Range("AW5").Value = 'The Sun River Nature Center experience increased our students'' understanding of birds of prey, desert habitat, interdependence of living things, animal adaptations, and space.  Fishing, nutrition, survival, and nature study sessions also helped students better understand science.'
$ws.Range("AX5").Value = 'During our Woodland Tour, students learn the importance of fire prevention, the impact of fires, the effect of litter on land and streams, soil erosion impact, laws and rules in the wilderness, and how to be a better environmental citizen.'
$ws.Range("AY5").Value = 'Before we leave, students  research and calculate the cost of Outdoor School including gas, meals, entrance fees, and camping fees.  They determine the total and cost per student.'
$ws.Range("AZ5").Value = 'Students articulate through songs and skits.  They also complete a Outdoor School journal and reflection that requires reading, writing, listening, and speaking skills.'
$ws.Range("BA5").Value = 'Students improve their social studies understanding through historical information and exhibits at the museum, nature center, and visitor centers.  Students learn about Native American culture, history, and crafts.'
$ws.Range("BB5").Value = 'All students, regardless of abilities, medical conditions, or socio-economic status can attend Outdoor School at no cost.  We work hard to meet these individual needs.'
$ws.Range("BD5").Value = 'For our limited English student, we provided a counselor who spoke Spanish and assisted the student.  Curriculum is translated as much as possible.'
$ws.Range("BE5").Value = 'Special education students receive assistance from their counselors.  Accommodations are provided to help them be successful. 1:1 if necessary'
$ws.Range("BF5").Value = 'Learning disabled students receive assistance from their counselors.  Accommodations are provided to help them be successful.  '
$ws.Range("BG5").Value = 'Items needed are purchased for the student.  Counselor helps with any needs.'
$ws.Range("BH5").Value = 'na'
$ws.Range("BI5").Value = 'na'
$ws.Range("BJ5").Value = 'na'
$ws.Range("BK5").Value = 'na'
$ws.Range("BL5").Value = 'Translator and translations when needed.  Spanish speaking counselor.'
$ws.Range("BM5").Value = 'Accommodations so they can fully participate.  (diabetes, asthma, epilepsy)'
$ws.Range("BN5").Value = 'At Amity Elementary School, all students are offered the opportunity to participate in Outdoor School.  Historically, we have only left behind students who are a danger to themselves or others and those with excessive missing work (even after numerous options and opportunities are given).  98% of 5th graders attend.'
$ws.Range("BP5").Value = 'Local fundraisers:  Stuff the Envelope, Community-School Dinners, Box-Tops, Jog-a-Thon, and Wreath Sales'

# Row 6
$ws.Range("A6").Value = 'Arlington School District #3'
$ws.Range("B6").Value = 'Arlington SD 3'
$ws.Range("C6").Value = 'Kevin Hunking'
$ws.Range("D6").Value = 'Kevin Hunking'
$ws.Range("E6").Value = 'khunking@arlington.k12.or.us'
$ws.Range("F6").Value = 'khunking@arlington.k12.or.us'
$ws.Range("G6").Value = '541-454-2632'
$ws.Range("H6").Value = '541-454-2632'
$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 5
$ws.Range("K6").Value = 5
$ws.Range("L6").Value = 3
$ws.Range("M6").Value = 4
$ws.Range("N6").Value = 4
$ws.Range("O6").Value = 2
$ws.Range("P6").Value = 3
$ws.Range("Q6").Value = 3
$ws.Range("R6").Value = 3
$ws.Range("S6").Value = 2
$ws.Range("T6").Value = 2
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = 4
$ws.Range("W6").Value = 5
$ws.Range("X6").Value = 1
$ws.Range("Y6").Value = 2
$ws.Range("AA6").Value = 'Follows core standards and integrates science, engineering, and math. Students learned orienteering and map reading.'
$ws.Range("AB6").Value = 'Teachers speak english and spanish.'
$ws.Range("AC6").Value = 'All students had equal access.'
$ws.Range("AD6").Value = 'Students learn about their communities and local environment.'
$ws.Range("AE6").Value = 'x'
$ws.Range("AF6").Value = 'x'
$ws.Range("AG6").Value = 'x'
$ws.Range("AH6").Value = 'x'
$ws.Range("AI6").Value = 'x'
$ws.Range("AJ6").Value = 'x'
$ws.Range("AK6").Value = 'x'
$ws.Range("AL6").Value = 'Enter "other" text here.'
$ws.Range("AM6").Value = 'Classroom discussion about outdoor school activities.'
$ws.Range("AN6").Value = 'Parents provided feedback and rangers helped teach classes.'
$ws.Range("AO6").Value = 'N/A'
$ws.Range("AP6").Value = 'Students received awards for completing achievements.'
$ws.Range("AQ6").Value = 'N/A'
$ws.Range("AR6").Value = 'Students received awards for completing achievements.'
$ws.Range("AS6").Value = 'Learned problem solving strategies as a group.'
$ws.Range("AT6").Value = 'N/A'
$ws.Range("AU6").Value = 'Learned to communicate with leaders and other students.'
$ws.Range("AV6").Value = 'N/A'
$ws.Range("AW6").Value = 'Learned orienteering and map reading.'
$ws.Range("AX6").Value = 'Learned to read maps.'
$ws.Range("AY6").Value = 'Map reading.'
$ws.Range("AZ6").Value = 'Reading and communication skills.'
$ws.Range("BA6").Value = 'Learned about our community and history.'
$ws.Range("BB6").Value = 'All students received differentiated instruction.'
$ws.Range("BC6").Value = 'Our students all live in a rural community.'
$ws.Range("BD6").Value = 'Instruction was available in other languages.'
$ws.Range("BE6").Value = 'Support staff was involved with special education students.'
$ws.Range("BF6").Value = 'Support staff was involved with special education students.'
$ws.Range("BG6").Value = 'All students received the same instruction.'
$ws.Range("BH6").Value = 'N/A'
$ws.Range("BI6").Value = 'N/A'
$ws.Range("BJ6").Value = 'N/A'
$ws.Range("BK6").Value = 'N/A'
$ws.Range("BL6").Value = 'Instruction was available in other languages.'
$ws.Range("BN6").Value = 'All students in grade 6 attend outdoor school as part of the curriculum.'
$ws.Range("BO6").Value = 'N/A'
$ws.Range("BP6").Value = 'Park Rangers donated instructional time.'

# Row 7
$ws.Range("A7").Value = 'AMS'
$ws.Range("B7").Value = 'Ashland SD 5'
$ws.Range("C7").Value = 'Hillary Harper'
$ws.Range("D7").Value = 'Alana Valencia'
$ws.Range("E7").Value = 'hillary.harper@ashland.k12.or.us'
$ws.Range("F7").Value = 'alana.valencia@ashland.k12.or.us'
$ws.Range("G7").Value = '541-482-1611'
$ws.Range("H7").Value = '541-481-2811'
$ws.Range("I7").Value = 5
$ws.Range("J7").Value = 1
$ws.Range("K7").Value = 4
$ws.Range("L7").Value = 4
$ws.Range("M7").Value = 5
$ws.Range("N7").Value = 5
$ws.Range("O7").Value = 3
$ws.Range("P7").Value = 1
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 4
$ws.Range("S7").Value = 2
$ws.Range("T7").Value = 3
$ws.Range("U7").Value = 3
$ws.Range("V7").Value = 4
$ws.Range("W7").Value = 3
$ws.Range("X7").Value = 2
$ws.Range("Y7").Value = 5
$ws.Range("AA7").Value = 'All lessons are alligned to Oregon State Standarsds and NGSS. '
$ws.Range("AB7").Value = 'Forms are produced in Spanish when necessary.'
$ws.Range("AC7").Value = 'All students are allowed to participate regardless of financial ability, scholarships and equipment are provided.'
$ws.Range("AD7").Value = 'n/a'
$ws.Range("AF7").Value = 'x'
$ws.Range("AH7").Value = 'x'
$ws.Range("AI7").Value = 'x'
$ws.Range("AJ7").Value = 'x'
$ws.Range("AL7").Value = 'Enter "other" text here.'
$ws.Range("AM7").Value = 'Pre and post assesments are given, comprehensive ecology unit is taught around the experience. Some teachers have received outdoor ed professional development regarding outdoor ed.  '
$ws.Range("AN7").Value = 'We have partnered with both the local high school to provide counselors and with the local university to provide instructors. Parents are encouraged to attend as volunteers counselors. '
$ws.Range("AO7").Value = 'n/a'
$ws.Range("AP7").Value = 'students come back from ROS with greater self sufficiency and leadership skills. '
$ws.Range("AQ7").Value = 'n/a'
$ws.Range("AR7").Value = 'increase in student engagement, students reflect fondly on the experience as they move through the middle school. '
$ws.Range("AS7").Value = 'we are able to work with student on positive problem solving skills which can be brought back into the classroom. '
$ws.Range("AT7").Value = 'students are able to creatively think about the interconnectedness of nature. '
$ws.Range("AU7").Value = 'living in a cabin with 8 other students over the course of the trip allows ample time to practice communication skills and ability to work. '
$ws.Range("AV7").Value = 'Students come back to school excited about their experience and this is carried through to all subject areas. '
$ws.Range("AW7").Value = 'NGSS and Oregon State Science Standards are taught.'
$ws.Range("AX7").Value = 'field studies allow students to practice real world science.'
$ws.Range("AY7").Value = 'some math  (measurement, angles, compass work, ) is practiced.'
$ws.Range("AZ7").Value = 'some writing and creative arts are used in field studies as well as presentation of findings. '
$ws.Range("BA7").Value = 'maps are used, and ways humans impact the environment are investigated.'
$ws.Range("BB7").Value = 'multiple intelligences and learning styles are addressed in field stations and the experience is available to all students regardless of ability. '
$ws.Range("BD7").Value = 'spanish forms are available'
$ws.Range("BE7").Value = 'EA support, 1 on 1 staffing, provided as necessary so all students can attend'
$ws.Range("BG7").Value = 'all students are allowed to attend regardless of financial ability, scholarships are available as well as equipment'
$ws.Range("BN7").Value = 'all students are allowed to attend regardless of financial ability, scholarships are available as well as equipment.EA support, 1 on 1 staffing, provided as necessary so all students can attend'
$ws.Range("BP7").Value = 'parent funding, district funding'

# Row 8
$ws.Range("A8").Value = 'Willow Wind'
$ws.Range("B8").Value = 'Ashland'
$ws.Range("C8").Value = 'Linda Terry'
$ws.Range("D8").Value = 'Alana Valencia'
$ws.Range("E8").Value = 'Linda.Terry@ashland.k12.or.us'
$ws.Range("F8").Value = 'alana.valencia@ashland.k12.or.us'
$ws.Range("G8").Value = '541-488-2684'
$ws.Range("H8").Value = '541-481-2811'
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 1
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 4
$ws.Range("M8").Value = 5
$ws.Range("N8").Value = 5
$ws.Range("O8").Value = 3
$ws.Range("P8").Value = 4
$ws.Range("Q8").Value = 3
$ws.Range("R8").Value = 5
$ws.Range("S8").Value = 2
$ws.Range("T8").Value = 2
$ws.Range("U8").Value = 3
$ws.Range("V8").Value = 5
$ws.Range("W8").Value = 5
$ws.Range("X8").Value = 2
$ws.Range("Y8").Value = 5
$ws.Range("AA8").Value = 'OMSI Outdoor Science School (the ODS used this year) provides a concept-based environmental program whereby activities are closely aligned with (and therefore are integrated with) Next Generation Science Standards, which have been implemented at Willow Wind.  This program takes place at Hancock Field Station, which is located on site at John Day Fossil Beds National Monument in Central Oregon.  Field studies include a day-long study of the geology/paleontology of John Day formations, and short-term courses will include avian studies, arid lands ecology, and orienteering, with evening programs in observational astronomy and the crepuscular habits of desert animals.  '
$ws.Range("AB8").Value = 'Several instructional leaders at OMSI Hancock Field Station are bilingual, so this option is presumably available if needed.'
$ws.Range("AC8").Value = 'Funding that covers the cost for all students to participate in ODS is the key to ensuring equity.  Without this support, many students would not be able to participate and reap the substantial benefits that ODS provides.'
$ws.Range("AD8").Value = 'Discussions between OMSI faciliators and students focused primarily on how urban areas were dependent on rural areas, with an emphasis on water use, food production and wildlands management.'
$ws.Range("AE8").Value = 'x'
$ws.Range("AF8").Value = 'x'
$ws.Range("AG8").Value = 'x'
$ws.Range("AH8").Value = 'x'
$ws.Range("AI8").Value = 'x'
$ws.Range("AJ8").Value = 'x'
$ws.Range("AK8").Value = 'x'
$ws.Range("AL8").Value = 'Enter "other" text here.'
$ws.Range("AM8").Value = 'Pre-activities consisted primarily of learning about John Day Fossil Beds and how students should prepare for their experience in a new environment.  Post-activites included extended studies of ecosystems dynamics, energy transfer in ecosystems and geologic processes/paleontology/geologic history.'
$ws.Range("AN8").Value = 'Parent volunteers assisted with organizing the ODS, and also volunteered as chaperones.'
$ws.Range("AO8").Value = 'Successful to some degree, although this is difficult to measure.'
$ws.Range("AP8").Value = 'Significant growth was noted among many students, particularly those who had not participated in overnight activites without parents.  Students learned that they could make decisions without input from parents, which built confidence.'
$ws.Range("AQ8").Value = 'Some students really blossomed in this environment-they were held to a high behavior expectation, and almost uniformly rose to the occasion.'
$ws.Range("AR8").Value = 'Successful, especially since students were expected to help with meal prep and cleanup; it was clear that for some students, this type of responsibility was new.'
$ws.Range("AS8").Value = 'Highly successful'
$ws.Range("AT8").Value = 'Highly successful'
$ws.Range("AU8").Value = 'Highly successful'
$ws.Range("AV8").Value = 'Greater enthusiasm for science'
$ws.Range("AW8").Value = 'Highly successful'
$ws.Range("AX8").Value = 'Highly successful'
$ws.Range("AY8").Value = 'n/a (minimal)'
$ws.Range("AZ8").Value = 'n/a'
$ws.Range("BA8").Value = 'n/a'
$ws.Range("BB8").Value = 'Somewhat successful'
$ws.Range("BG8").Value = 'Trip equipment and supplies were provided to students who required this support.'
$ws.Range("BN8").Value = 'Support will continue to be provided so that all students at Willow Wind can participate in ODS.'
$ws.Range("BP8").Value = 'Ashland School District covered any costs not covered through grant monies.  Supplies for students were provided by parents and school staff.'

# Row 9
$ws.Range("A9").Value = 'John Muir'
$ws.Range("B9").Value = 'Ashland SD 5'
$ws.Range("C9").Value = 'James Bowers   Marcia Ososke'
$ws.Range("D9").Value = 'Alana Valencia'
$ws.Range("E9").Value = 'James.Bowers@Ashland.k12.or.us    Marcia.Ososke@Ashland.k12.or.us'
$ws.Range("F9").Value = 'alana.valencia@ashland.k12.or.us'
$ws.Range("G9").Value = '541-482-8577'
$ws.Range("H9").Value = '541-481-2811'
$ws.Range("I9").Value = 5
$ws.Range("J9").Value = 5
$ws.Range("K9").Value = 5
$ws.Range("L9").Value = 5
$ws.Range("M9").Value = 5
$ws.Range("N9").Value = 5
$ws.Range("O9").Value = 4
$ws.Range("P9").Value = 2
$ws.Range("Q9").Value = 5
$ws.Range("R9").Value = 5
$ws.Range("S9").Value = 3
$ws.Range("T9").Value = 4
$ws.Range("U9").Value = 4
$ws.Range("V9").Value = 4
$ws.Range("W9").Value = 4
$ws.Range("X9").Value = 2
$ws.Range("Y9").Value = 4
$ws.Range("Z9").Value = 3
$ws.Range("AA9").Value = 'Next Generation Science Standards:  MS-LS2-4;  MS-LS2-5;  MS-LS4-2; MS-ESS2-4; MS-ESS3-3;  MS-Ess3-4;  MS-ESS3-5'
$ws.Range("AB9").Value = 'Not offered bilingually.'
$ws.Range("AC9").Value = 'We provided students with the opportunity to study "Ecosystems" that were distant from our region.  In traveling to the Oregon Coast, we studied the differences in ecosystems in regards to forests, water, ocean and estuaries.  Students learned not only about "what makes an ecosystem," but about current events that effect the health of ecosystems from the coast to the Rogue Valley.  Students visited the coast, hiking along part of the Oregon Coast Trail, visiting Shore Acres State Park, Cape Arago, Sunset Bay State Park and South Slough National Marine Estuary.'
$ws.Range("AD9").Value = 'See attached.  Students worked with focus questions looking at how rural economies are dependent upon the natural environment.  We invited biologists from Oregon Department of Fish and Wildlife to discuss work on rural watersheds to manage for habitat while sustaining the local environment.  At South Slough Estuary, students learned about local costal watersheds and the factors that impact them in regards to maintaining biodiversity.  '
$ws.Range("AE9").Value = 'x'
$ws.Range("AF9").Value = 'x'
$ws.Range("AG9").Value = 'x'
$ws.Range("AH9").Value = 'x'
$ws.Range("AI9").Value = 'x'
$ws.Range("AJ9").Value = 'x'
$ws.Range("AK9").Value = 'x'
$ws.Range("AL9").Value = 'Hands-on experiences through hiking, observing, journaling and discussing.'
$ws.Range("AM9").Value = 'Students completed experiential education reflections regarding their involvement in the trip, along with observations and big-picture learnings.  '
$ws.Range("AN9").Value = 'Parents joined on this trip to assist as chaperones in yurts and tents and to help with managing students, along with preparing meals.  We worked with Sunset Bay State Park rangers to do a service learning project, removing blackberries around the park, as a means of aquiring a discount on rental rates.  We also made the connection between managing invasive species in an ecosystem.  We coordinated with educators at South Slough National Marine Estuary in guided hikes to learn about the importance of the estuary ecosystem in the region.  Rangers from the Oregon Department of Fish and Wildlife met with our group one evening to teach about ecosystem management, we did a tour of the Washed Ashore museum in Bandon to learn about human impacts on the ocean environment on our last day.  Also, on the third day we met with an educator from Costal Ecosystem Partners in Education to learn about bird migration and habitat protection.  On our first night, we took students to visit Siuslaw Tribal Elders in Coos Bay to learn about Native Peoples'' relationship to the land and the transition of land ownership in the area.  Students learned native ways of agriculture, hunting and ceremonial traditions.'
$ws.Range("AO9").Value = 'As a school we score well on district assessments (Easy CBM) and state/national assessements (Smarter Balanced)'
$ws.Range("AP9").Value = 'Doing a fall over night trip in particular helps build leadership skills right from the start of the school year.  Every year we see students brave the elements and become more confident people.'
$ws.Range("AQ9").Value = 'We see classmates bond and not want to let each other down this leads to fewer discipline and management problems.'
$ws.Range("AR9").Value = 'I think so but no data (quantitive or quailitative) around this.  If there was a survey of something for the students to fill out after the trip it may give us this type of information. '
$ws.Range("AS9").Value = 'I think so but no data (quantitive or quailitative) around this.  If there was a survey of something for the students to fill out after the trip it may give us this type of information. '
$ws.Range("AT9").Value = 'I think so but no data (quantitive or quailitative) around this.  If there was a survey of something for the students to fill out after the trip it may give us this type of information. '
$ws.Range("AU9").Value = 'I think so but no data (quantitive or quailitative) around this.  If there was a survey of something for the students to fill out after the trip it may give us this type of information. '
$ws.Range("AV9").Value = 'I think so but no data (quantitive or quailitative) around this.  If there was a survey of something for the students to fill out after the trip it may give us this type of information. '
$ws.Range("AW9").Value = 'I think so but no data (quantitive or quailitative) around this.  If there was a survey of something for the students to fill out after the trip it may give us this type of information. '
$ws.Range("AX9").Value = 'I think so but no data (quantitive or quailitative) around this.  If there was a survey of something for the students to fill out after the trip it may give us this type of information. '
$ws.Range("AY9").Value = 'I think so but no data (quantitive or quailitative) around this.  If there was a survey of something for the students to fill out after the trip it may give us this type of information. '
$ws.Range("AZ9").Value = 'I think so but no data (quantitive or quailitative) around this.  If there was a survey of something for the students to fill out after the trip it may give us this type of information. '
$ws.Range("BA9").Value = 'I think so but no data (quantitive or quailitative) around this.  If there was a survey of something for the students to fill out after the trip it may give us this type of information. '
$ws.Range("BB9").Value = 'Outdoor expereinces all allow for all students to thrive.  This especially goes for students who may struggle more through tradtional acedemic classroom expereinces.'
$ws.Range("BD9").Value = 'Visuals and hands-on experiences.  Lots of opportunities to work in groups.'
$ws.Range("BE9").Value = 'Educational assistants joined on trip to assist with learnings.  Activities were modified for students in terms of pace and directions.  '
$ws.Range("BF9").Value = 'No major disabilities with students - lots of parent support and staff support for those who might have struggled either emotionally or socially.'
$ws.Range("BG9").Value = 'Gear provided by the school for all students.  Food provided.'
$ws.Range("BH9").Value = 'Connected students to Native American tribal traditions with Siuslaw Tribe in Coos Bay.'
$ws.Range("BI9").Value = 'N/A'
$ws.Range("BJ9").Value = 'N/A'
$ws.Range("BK9").Value = 'N/A'
$ws.Range("BL9").Value = 'N/A'
$ws.Range("BN9").Value = 'This trip is required for all students in our school.  We provide gear, food, lodging and transportation for all students, regardless of economic level.  We fundraise as part of our Parent Teacher Collective to scholarship any needs students or families may have.  It is so great to have the outdoor school funds as it takes some financial pressure off of the district and our own Parent Teacher Collective.  Curriculum in the classroom was integrated with trip studies.'
$ws.Range("BO9").Value = 'N/A'
$ws.Range("BP9").Value = 'Parents attended trip, using personal vehicles and donating gas costs.  Cooking equipment and gear were provided by parents.  Coordinated with Sunset Bay State Park to do community service for discounted camping rates.  District pays for fuel for busses.'

